# "feat: enable module API docs: work on rapport"
# Update the "Journal de travail" time log: two more entries worked on the
# report (Rapport) bring C30/C31 up to 3h each, and the now-unused helper
# rows/formatting below the table (C32:C38, left over from an old fill-down)
# are cleared out.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Journal de travail")
$ws.Activate()

# The leftover stray formatting below the table (C30:C38) is cleared out
# entirely first (contents + number format) ...
$ws.Range("C30:C38").Clear()

# ... then the two real entries get their updated hours.
# Row 30 (Implémentation / Rapport): 2h -> 3h
$ws.Range("C30").Value = 3
# Row 31 (Rédaction / Rapport): was blank -> 3h
$ws.Range("C31").Value = 3

# Leave the selection where work left off.
$ws.Range("D33").Select()
